# Insert a new data row at row 706 on the single worksheet, shifting the
# existing rows 706-796 down to 707-797 (dimension grows from A1:R796 to
# A1:R797), then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 706..796 down by one, creating a blank row at 706.
$ws.Rows(706).Insert()

# Populate the newly inserted row 706 with the new record's data.
$ws.Range("A706").Value = 3
$ws.Range("B706").Value = "Femacal de La Calera"
$ws.Range("C706").Value = "Coquimbo"
$ws.Range("D706").Value = 45142
$ws.Range("E706").Value = 5
$ws.Range("F706").Value = 100112037
$ws.Range("G706").Value = "Cebollín"
$ws.Range("H706").Value = "Sin especificar"
$ws.Range("I706").Value = "Primera"
$ws.Range("J706").Value = 190
$ws.Range("K706").Value = 4000
$ws.Range("L706").Value = 4300
$ws.Range("M706").Value = 4111
$ws.Range("N706").Value = "$/paquete 36 unidades"
$ws.Range("O706").Value = "Provincia de Quillota"
$ws.Range("P706").Value = 114
$ws.Range("Q706").Value = 36
$ws.Range("R706").Value = "Hortaliza"

# Keep the date column formatted like the rest of column D.
$ws.Range("D706").NumberFormat = $ws.Range("D707").NumberFormat
